$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (English)
$ws.Range("B2").Value = "Vostro"
$ws.Range("C2").Value = "Dell"
$ws.Range("D2").Value = 3568
$ws.Range("E2").Value = "DKS"
$ws.Range("F2").Value = 1.454
$ws.Range("G2").Value = "To take enrollments"

# Row 3 (Arabic)
$ws.Range("B3").Value = "ستر  "
$ws.Range("C3").Value = "دلّ  "
$ws.Range("D3").Value = 3568
$ws.Range("E3").Value = "DKS"
$ws.Range("F3").Value = 1.454
$ws.Range("G3").Value = "لأخذ التسجيلات"

# View settings: scroll to column C and select from A4 downward
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("A4:A1048576").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
